# The post "「気分が良いぞう」I فيل GOOD" (row 858) was removed from the blog.
# Delete its entire row; Excel shifts every following row up by one,
# which also shrinks the sheet's used range from A1:C869 to A1:C868.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(858).Delete()
